$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily disk-usage samples appended after the existing data (rows 2-106).
# Columns: A=date (serial), B=time (fraction of day), C=files, D=disk_space
$newRows = @(
    @{ Row = 107; Date = 45126; Time = 0.68236111111111108; Files = 81639; Disk = 1630 },
    @{ Row = 108; Date = 45127; Time = 0.45068287037037041; Files = 83211; Disk = 1680 },
    @{ Row = 109; Date = 45128; Time = 0.48158564814814814; Files = 83215; Disk = 1680 },
    @{ Row = 110; Date = 45129; Time = 0.58716435185185178; Files = 83215; Disk = 1680 },
    @{ Row = 111; Date = 45130; Time = 0.47083333333333338; Files = 83215; Disk = 1680 },
    @{ Row = 112; Date = 45131; Time = 0.51880787037037035; Files = 83217; Disk = 1680 },
    @{ Row = 113; Date = 45132; Time = 0.60712962962962969; Files = 83217; Disk = 1680 }
)

# The last existing data row carries the date/time number formats (m/d/yyyy and
# h:mm:ss) that every row in the table uses. Copy just those formats down onto
# each new row before writing values, so the new cells reuse the same styles
# instead of minting new ones.
$fmtSource = $ws.Range("A106:D106")

foreach ($r in $newRows) {
    $fmtSource.Copy()
    $ws.Range("A" + $r.Row + ":D" + $r.Row).PasteSpecial(-4122)

    $ws.Range("A" + $r.Row).Value = $r.Date
    $ws.Range("B" + $r.Row).Value = $r.Time
    $ws.Range("C" + $r.Row).Value = $r.Files
    $ws.Range("D" + $r.Row).Value = $r.Disk
}

$excel.CutCopyMode = $false

# Match the author's view state: scrolled so row 93 is at the top, with the
# next empty row's first cell selected.
$excel.ActiveWindow.ScrollRow = 93
$ws.Range("A114").Select()
